$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# The "A Suite" row in the module table is now the "IAM" row.
$ws.Range("A2").Value = "IAM"

# Update the active selection to match the author's saved view (C3 instead of C4)
$ws.Activate()
$ws.Range("C3").Select()
